$d = $word.ActiveDocument

function Replace-WithOoxml($findText, $bodyXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $findText"
    }

    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $bodyXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($pkg)
}

# --- Title: "Communicatie-log externe- en interne ontwikkelaar" ---
#     -> "Communicationlog external- and internal developer" (spell-checked, split into words)
$titleXml =
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Communicationlog</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>external</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>and</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>internal</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>developer</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'

Replace-WithOoxml "Communicatie-log externe- en interne ontwikkelaar" $titleXml

# --- Table header: "Wanneer?" -> "When?" ---
$whenXml =
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>When</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>?</w:t></w:r>'

Replace-WithOoxml "Wanneer?" $whenXml

# --- Table header: "Hoe lang?" -> "How long?" ---
$howLongXml =
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>How long</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>?</w:t></w:r>'

Replace-WithOoxml "Hoe lang?" $howLongXml

# --- Table header: "Wat is besproken?" -> "What?" ---
$whatXml =
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>What</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>?</w:t></w:r>'

Replace-WithOoxml "Wat is besproken?" $whatXml
